$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dimension")
$ws2 = $wb.Worksheets.Item("Metrics")

# --- Sheet2 (Metrics): remove the "dfa:activeViewViewableImpressions" row (row 1) ---
$ws2.Rows.Item(1).Delete()

# --- Sheet1 (Dimension): swap the last two rows (A18 <-> A19) ---
$valA18 = $ws1.Range("A18").Value()
$valA19 = $ws1.Range("A19").Value()
$ws1.Range("A18").Value = $valA19
$ws1.Range("A19").Value = $valA18

# --- Update selections / active sheet / active cell ---
$ws2.Range("A1:XFD1").Select()
$ws1.Activate()
$ws1.Range("A18").Select()

$wb.Save()
